$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so that numeric-looking
# strings (e.g. "1.001", "23.219.61") are preserved exactly as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '23.219.61'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '1.606.88'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '0.9998'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').Value = '303.70'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('D7').Value = '0.3773'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = '52.11'
$ws.Range('E8').Value = '  +4.53%  '
$ws.Range('D9').Value = '0.3638'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').Value = '1.275'
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('D11').Value = '0.08141'
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = '0.9992'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '22.82'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').Value = '6.589'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '7.404'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '0.00001251'
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('D17').Value = '1.608.71'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '94.06'
$ws.Range('E18').Value = '  +2.26%  '
$ws.Range('D19').Value = '0.06936'
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').Value = '18.16'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').Value = '6.540'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').Value = '12.93'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').Value = '23.223.80'
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('D25').Value = '2.447'
$ws.Range('E25').Value = '  +4.07%  '
$ws.Range('D26').Value = '3.057'
$ws.Range('E26').Value = '  +8.66%  '
$ws.Range('D27').Value = '21.22'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').Value = '149.72'
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('D29').Value = '5.285'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('D30').Value = '135.46'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').Value = '2.379'
$ws.Range('E31').Value = '  +2.87%  '
$ws.Range('D32').Value = '6.763'
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('D33').Value = '1.778.20'
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('D34').Value = '0.9651'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').Value = '0.07500'
$ws.Range('E35').Value = '  -1.19%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '0.02763'
$ws.Range('E36').Value = '  +1.76%  '
$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').Value = '10.34'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').Value = '0.2521'
$ws.Range('E38').Value = '  -0.63%  '
$ws.Range('D39').Value = '6.143'
$ws.Range('E39').Value = '  -2.35%  '
$ws.Range('D40').Value = '0.08804'
$ws.Range('E40').Value = '  -1.02%  '
$ws.Range('D41').Value = '1.389'
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('D42').Value = '0.7098'
$ws.Range('E42').Value = '  +0.77%  '
$ws.Range('D43').Value = '12.49'
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('D44').Value = '15.79'
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('D45').Value = '0.6541'
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('D46').Value = '2.333'
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('D47').Value = '0.9986'
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('D48').Value = '4.007'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').Value = '132.87'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').Value = '0.07947'
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('D51').Value = '1.208'
$ws.Range('E51').Value = '  -1.98%  '

# Restore original (default/general) formatting so the cell style
# matches the unedited cells (no explicit style index).
$ws.Range("D2:D51").ClearFormats()
